$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 398.92856
$ws.Range("I80").Value = 308.33334
$ws.Range("J80").Value = 423.63635
$ws.Range("K80").Value = 925.0000200000001
$ws.Range("L80").Value = 1270.90905
$ws.Range("M80").Value = 72.99997999999994
$ws.Range("N80").Value = -3266.90905
$ws.Range("H83").Value = 398.92856
$ws.Range("I83").Value = 308.33334
$ws.Range("J83").Value = 423.63635
$ws.Range("K83").Value = 2775.00006
$ws.Range("L83").Value = 3812.72715
$ws.Range("M83").Value = 2216.99994
$ws.Range("N83").Value = -13796.72715
$ws.Range("H88").Value = 1500
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 1500
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 1500
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -2312
$ws.Range("H91").Value = 1500
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 1500
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 1500
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -4308
$ws.Range("H96").Value = 961.3333
$ws.Range("I96").Value = 1292.25
$ws.Range("J96").Value = 299.5
$ws.Range("K96").Value = 3876.75
$ws.Range("L96").Value = 898.5
$ws.Range("M96").Value = -2503.75
$ws.Range("N96").Value = -3644.5
$ws.Range("H112").Value = 3199.1667
$ws.Range("I112").Value = 1015.5
$ws.Range("J112").Value = 3563.111
$ws.Range("K112").Value = 3046.5
$ws.Range("L112").Value = 10689.333
$ws.Range("M112").Value = -1938.5
$ws.Range("N112").Value = -12905.333
$ws.Range("H113").Value = 15000.182
$ws.Range("J113").Value = 15625.25
$ws.Range("L113").Value = 15625.25
$ws.Range("N113").Value = -22133.25
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").ClearContents()
$ws.Range("H132").Value = 73800.71000000001
$ws.Range("I132").Value = 2029.2858
$ws.Range("J132").Value = 145572.14
$ws.Range("K132").Value = 6087.857400000001
$ws.Range("L132").Value = 436716.42
$ws.Range("M132").Value = -3557.857400000001
$ws.Range("N132").Value = -441776.42
$ws.Range("H135").Value = 1329.0869
$ws.Range("I135").Value = 937.82355
$ws.Range("K135").Value = 8440.41195
$ws.Range("M135").Value = -5905.41195
$ws.Range("H137").Value = 1685.738
$ws.Range("I137").Value = 1568.1154
$ws.Range("K137").Value = 4704.3462
$ws.Range("M137").Value = -2154.3462

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 5999.3335
$ws.Range("J25").Value = 6999.5
$ws.Range("L25").Value = 6999.5
$ws.Range("N25").Value = -7803.5
$ws.Range("H61").Value = 2786.524
$ws.Range("I61").Value = 2515.7354
$ws.Range("J61").Value = 3937.375
$ws.Range("K61").Value = 2515.7354
$ws.Range("L61").Value = 3937.375
$ws.Range("M61").Value = -2303.7354
$ws.Range("N61").Value = -4361.375
$ws.Range("H74").Value = 5989.8335
$ws.Range("I74").Value = 1601.7916
$ws.Range("K74").Value = 1601.7916
$ws.Range("M74").Value = -727.7916
$ws.Range("H77").Value = 5989.8335
$ws.Range("I77").Value = 1601.7916
$ws.Range("K77").Value = 8008.958000000001
$ws.Range("M77").Value = -3640.958000000001
$ws.Range("H132").Value = 2529.4146
$ws.Range("I132").Value = 2263.0293
$ws.Range("J132").Value = 3823.2856
$ws.Range("K132").Value = 6789.0879
$ws.Range("L132").Value = 11469.8568
$ws.Range("M132").Value = -4259.0879
$ws.Range("N132").Value = -16529.8568
$ws.Range("H136").Value = 2786.524
$ws.Range("I136").Value = 2515.7354
$ws.Range("J136").Value = 3937.375
$ws.Range("K136").Value = 7547.206200000001
$ws.Range("L136").Value = 11812.125
$ws.Range("M136").Value = -4997.206200000001
$ws.Range("N136").Value = -16912.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3387.2554
$ws.Range("I20").Value = 2967.9
$ws.Range("J20").Value = 4127.294
$ws.Range("K20").Value = 2967.9
$ws.Range("L20").Value = 4127.294
$ws.Range("M20").Value = -2720.9
$ws.Range("N20").Value = -4621.294
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H86").Value = 2393.6177
$ws.Range("I86").Value = 2366.52
$ws.Range("J86").Value = 2468.889
$ws.Range("K86").Value = 2366.52
$ws.Range("L86").Value = 2468.889
$ws.Range("M86").Value = -1243.52
$ws.Range("N86").Value = -4714.889
$ws.Range("H89").Value = 2393.6177
$ws.Range("I89").Value = 2366.52
$ws.Range("J89").Value = 2468.889
$ws.Range("K89").Value = 11832.6
$ws.Range("L89").Value = 12344.445
$ws.Range("M89").Value = -6216.6
$ws.Range("N89").Value = -23576.445
$ws.Range("H99").Value = 3846.2173
$ws.Range("I99").Value = 1970.1666
$ws.Range("J99").Value = 5892.8184
$ws.Range("K99").Value = 1970.1666
$ws.Range("L99").Value = 5892.8184
$ws.Range("M99").Value = -472.1666
$ws.Range("N99").Value = -8888.8184
$ws.Range("H107").Value = 885.70966
$ws.Range("I107").Value = 799.38464
$ws.Range("J107").Value = 1334.6
$ws.Range("K107").Value = 799.38464
$ws.Range("L107").Value = 1334.6
$ws.Range("M107").Value = 1120.61536
$ws.Range("N107").Value = -5174.6
$ws.Range("H134").Value = 2169.111
$ws.Range("I134").Value = 2162.3508
$ws.Range("K134").Value = 6487.0524
$ws.Range("M134").Value = -3952.0524

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 103438.1
$ws.Range("I31").Value = 127576.5
$ws.Range("J31").Value = 6884.5
$ws.Range("K31").Value = 127576.5
$ws.Range("L31").Value = 6884.5
$ws.Range("M31").Value = -127281.5
$ws.Range("N31").Value = -7474.5
$ws.Range("H34").Value = 103438.1
$ws.Range("I34").Value = 127576.5
$ws.Range("J34").Value = 6884.5
$ws.Range("K34").Value = 127576.5
$ws.Range("L34").Value = 6884.5
$ws.Range("M34").Value = -127374.5
$ws.Range("N34").Value = -7288.5
$ws.Range("H62").Value = 7905.1113
$ws.Range("J62").Value = 8451.4
$ws.Range("L62").Value = 8451.4
$ws.Range("N62").Value = -9699.4
$ws.Range("H65").Value = 7905.1113
$ws.Range("J65").Value = 8451.4
$ws.Range("L65").Value = 42257
$ws.Range("N65").Value = -48497
$ws.Range("H68").Value = 100000
$ws.Range("I68").Value = 100000
$ws.Range("K68").Value = 100000
$ws.Range("M68").Value = -99251
$ws.Range("H71").Value = 100000
$ws.Range("I71").Value = 100000
$ws.Range("K71").Value = 300000
$ws.Range("M71").Value = -296256
$ws.Range("H88").Value = 22500
$ws.Range("I88").Value = 10000
$ws.Range("J88").Value = 26666.666
$ws.Range("K88").Value = 10000
$ws.Range("L88").Value = 26666.666
$ws.Range("M88").Value = -9594
$ws.Range("N88").Value = -27478.666
$ws.Range("H91").Value = 22500
$ws.Range("I91").Value = 10000
$ws.Range("J91").Value = 26666.666
$ws.Range("K91").Value = 10000
$ws.Range("L91").Value = 26666.666
$ws.Range("M91").Value = -8596
$ws.Range("N91").Value = -29474.666
$ws.Range("H107").Value = 878.2
$ws.Range("I107").Value = 618.75
$ws.Range("K107").Value = 618.75
$ws.Range("M107").Value = 1301.25
$ws.Range("H132").Value = 2986.9678
$ws.Range("I132").Value = 2993.4285
$ws.Range("J132").Value = 2973.4
$ws.Range("K132").Value = 8980.2855
$ws.Range("L132").Value = 8920.200000000001
$ws.Range("M132").Value = -6450.2855
$ws.Range("N132").Value = -13980.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H124").Value = 32255.334
$ws.Range("J124").Value = 40166.668
$ws.Range("L124").Value = 120500.004
$ws.Range("N124").Value = -130320.004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1788.8966
$ws.Range("I113").Value = 2168.762
$ws.Range("J113").Value = 791.75
$ws.Range("K113").Value = 2168.762
$ws.Range("L113").Value = 791.75
$ws.Range("M113").Value = 1.237999999999829
$ws.Range("N113").Value = -5131.75
$ws.Range("H122").Value = 2394.5
$ws.Range("I122").Value = 2140
$ws.Range("K122").Value = 6420
$ws.Range("M122").Value = -3970
$ws.Range("H126").Value = 10843.462
$ws.Range("I126").Value = 13205.15
$ws.Range("J126").Value = 2971.1667
$ws.Range("K126").Value = 39615.45
$ws.Range("L126").Value = 8913.500100000001
$ws.Range("M126").Value = -37145.45
$ws.Range("N126").Value = -13853.5001
$ws.Range("H136").Value = 30511.732
$ws.Range("J136").Value = 30511.732
$ws.Range("L136").Value = 91535.196
$ws.Range("N136").Value = -96635.196

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 31916.666
$ws.Range("J63").Value = 32875
$ws.Range("L63").Value = 32875
$ws.Range("N63").Value = -34373
$ws.Range("H66").Value = 31916.666
$ws.Range("J66").Value = 32875
$ws.Range("L66").Value = 98625
$ws.Range("N66").Value = -106113
$ws.Range("H68").Value = 3843.375
$ws.Range("I68").Value = 2791.6667
$ws.Range("J68").Value = 6998.5
$ws.Range("K68").Value = 2791.6667
$ws.Range("L68").Value = 6998.5
$ws.Range("M68").Value = -2042.6667
$ws.Range("N68").Value = -8496.5
$ws.Range("H71").Value = 3843.375
$ws.Range("I71").Value = 2791.6667
$ws.Range("J71").Value = 6998.5
$ws.Range("K71").Value = 13958.3335
$ws.Range("L71").Value = 34992.5
$ws.Range("M71").Value = -10214.3335
$ws.Range("N71").Value = -42480.5
$ws.Range("H136").Value = 4076.375
$ws.Range("I136").Value = 3786.3076
$ws.Range("K136").Value = 11358.9228
$ws.Range("M136").Value = -8808.9228

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1637.7858
$ws.Range("I132").Value = 1636.6285
$ws.Range("J132").Value = 1643.5714
$ws.Range("K132").Value = 4909.8855
$ws.Range("L132").Value = 4930.7142
$ws.Range("M132").Value = -2379.8855
$ws.Range("N132").Value = -9990.7142
$ws.Range("H136").Value = 3215.3845
$ws.Range("I136").Value = 3264.4443
$ws.Range("J136").Value = 3105
$ws.Range("K136").Value = 9793.332900000001
$ws.Range("L136").Value = 9315
$ws.Range("M136").Value = -7243.332900000001
$ws.Range("N136").Value = -14415
